# Redid prov indicators; added elasticity analysis
#
# The workbook has 9 sheets (count, null_count, mean, std, min, 25%, 50%,
# 75%, max), each backed by its own Excel Table whose column 28 (AB) is the
# provider indicator column "Provider_PPSA". This renames that column to
# "Provider_MPSA" on every sheet (which also updates the shared table
# column metadata since the header cell drives the ListColumn name),
# widens column AB slightly everywhere to fit the new header text, and
# refreshes the recomputed mean/std statistics (sheet "mean" / "std") plus
# the 50%/75% quantile rows to reflect the re-derived indicator.

$wb = $excel.ActiveWorkbook

$newHeader = "Provider_MPSA"

# ---------------------------------------------------------------------
# 1) Rename the "Provider_PPSA" column header to "Provider_MPSA" on every
#    sheet's table (column AB = column 28), then fix up column AB's width.
# ---------------------------------------------------------------------
for ($i = 1; $i -le 9; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $ws.Cells.Item(1, 28).Value = $newHeader
}

# Column-AB widths after the rename: every sheet goes from its old bestFit
# width to a value one unit wider, except the "std" sheet which was
# already at the wider value.
$ws1 = $wb.Worksheets.Item(1)   # count
$ws1.Columns.Item(28).ColumnWidth = 16.166666666666668   # -> stored width 17

$ws2 = $wb.Worksheets.Item(2)   # null_count
$ws2.Columns.Item(28).ColumnWidth = 16.166666666666668   # -> stored width 17

$ws3 = $wb.Worksheets.Item(3)   # mean
$ws3.Columns.Item(28).ColumnWidth = 19.166666666666668   # -> stored width 20

$ws5 = $wb.Worksheets.Item(5)   # min
$ws5.Columns.Item(28).ColumnWidth = 16.166666666666668   # -> stored width 17

$ws6 = $wb.Worksheets.Item(6)   # 25%
$ws6.Columns.Item(28).ColumnWidth = 16.166666666666668   # -> stored width 17

$ws7 = $wb.Worksheets.Item(7)   # 50%
$ws7.Columns.Item(28).ColumnWidth = 16.166666666666668   # -> stored width 17

$ws8 = $wb.Worksheets.Item(8)   # 75%
$ws8.Columns.Item(28).ColumnWidth = 16.166666666666668   # -> stored width 17

$ws9 = $wb.Worksheets.Item(9)   # max
$ws9.Columns.Item(28).ColumnWidth = 16.166666666666668   # -> stored width 17

# sheet 4 ("std") keeps its existing width (already the wider value) - no
# change needed there.

# ---------------------------------------------------------------------
# 2) Recomputed statistics for column AB (rows 2-19) on the sheets whose
#    numbers actually move once the indicator was redefined: mean, std,
#    50% and 75%. (count, null_count, min, 25% and max are unaffected.)
# ---------------------------------------------------------------------

$meanValues = @(
    0.116504854368932, 0.116504854368932, 0.116504854368932,
    0.116504854368932, 0.116504854368932, 0.116504854368932,
    0.1176470588235294, 0.1176470588235294,
    0.1188118811881188, 0.1188118811881188, 0.1188118811881188, 0.1188118811881188,
    0.12, 0.12,
    0.1212121212121212, 0.1212121212121212, 0.1212121212121212,
    0.1263157894736842
)

$stdValues = @(
    0.3223982093132047, 0.3223982093132047, 0.3223982093132047,
    0.3223982093132047, 0.3223982093132047, 0.3223982093132048,
    0.3237808098282633, 0.3237808098282633,
    0.3251808331642962, 0.3251808331642962, 0.3251808331642962, 0.3251808331642962,
    0.3265986323710905, 0.3265986323710905,
    0.328034569878314, 0.328034569878314, 0.328034569878314,
    0.333967295607331
)

for ($r = 2; $r -le 19; $r++) {
    $ws3.Cells.Item($r, 28).Value = $meanValues[$r - 2]
}

for ($r = 2; $r -le 19; $r++) {
    $ws4 = $wb.Worksheets.Item(4)   # std
    $ws4.Cells.Item($r, 28).Value = $stdValues[$r - 2]
}

for ($r = 2; $r -le 19; $r++) {
    $ws7.Cells.Item($r, 28).Value = 0
}

for ($r = 2; $r -le 19; $r++) {
    $ws8.Cells.Item($r, 28).Value = 0
}
